$wb = $excel.ActiveWorkbook

# Update the "% Completed" column for the Database tracking sheet
$wsDatabase = $wb.Worksheets.Item("Database")
$wsDatabase.Range("D5").Value = 1
$wsDatabase.Range("D6").Value = 1

# Leave the cursor where the editor last left it on the Document sheet
$wsDocument = $wb.Worksheets.Item("Document")
$wsDocument.Activate() | Out-Null
$wsDocument.Range("B27").Select() | Out-Null

# Finish on the Database sheet, which becomes the active tab
$wsDatabase.Activate() | Out-Null
$wsDatabase.Range("G7").Select() | Out-Null
